# Update the "想去人数" (number of people interested) figures that changed
# between the two data pulls, on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 6047
$wsExhibit.Range("F9").Value = 35
$wsExhibit.Range("F10").Value = 62
$wsExhibit.Range("F12").Value = 140
$wsExhibit.Range("F13").Value = 340
$wsExhibit.Range("F14").Value = 440
$wsExhibit.Range("F15").Value = 3048
$wsExhibit.Range("F17").Value = 171
$wsExhibit.Range("F18").Value = 1680
$wsExhibit.Range("F19").Value = 16

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 6047
$wsAll.Range("F10").Value = 35
$wsAll.Range("F11").Value = 62
$wsAll.Range("F13").Value = 140
$wsAll.Range("F14").Value = 340
$wsAll.Range("F15").Value = 440
$wsAll.Range("F16").Value = 3048
$wsAll.Range("F18").Value = 171
$wsAll.Range("F19").Value = 1680
$wsAll.Range("F20").Value = 16

$wb.Save()
